$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$new = '44 x 11' + [char]11 + '  1    1' + [char]11 + '  ----' + [char]11 + '4|    |' + [char]11 + '4|    |'
$cell.Range.Text = $new

$cell = $t.Cell(1, 2)
$new = '48 x 34' + [char]11 + '  3    4' + [char]11 + '  ----' + [char]11 + '4|    |' + [char]11 + '8|    |'
$cell.Range.Text = $new

$cell = $t.Cell(1, 3)
$new = '66 x 68' + [char]11 + '  6    8' + [char]11 + '  ----' + [char]11 + '6|    |' + [char]11 + '6|    |'
$cell.Range.Text = $new

$cell = $t.Cell(2, 1)
$new = '68 x 74' + [char]11 + '  7    4' + [char]11 + '  ----' + [char]11 + '6|    |' + [char]11 + '8|    |'
$cell.Range.Text = $new

$cell = $t.Cell(2, 2)
$new = '64 x 34' + [char]11 + '  3    4' + [char]11 + '  ----' + [char]11 + '6|    |' + [char]11 + '4|    |'
$cell.Range.Text = $new

$cell = $t.Cell(2, 3)
$new = '38 x 77' + [char]11 + '  7    7' + [char]11 + '  ----' + [char]11 + '3|    |' + [char]11 + '8|    |'
$cell.Range.Text = $new

$cell = $t.Cell(3, 1)
$new = '11 x 58' + [char]11 + '  5    8' + [char]11 + '  ----' + [char]11 + '1|    |' + [char]11 + '1|    |'
$cell.Range.Text = $new

$cell = $t.Cell(3, 2)
$new = '62 x 70' + [char]11 + '  7    0' + [char]11 + '  ----' + [char]11 + '6|    |' + [char]11 + '2|    |'
$cell.Range.Text = $new

$cell = $t.Cell(3, 3)
$new = '53 x 18' + [char]11 + '  1    8' + [char]11 + '  ----' + [char]11 + '5|    |' + [char]11 + '3|    |'
$cell.Range.Text = $new

$cell = $t.Cell(4, 1)
$new = '55 x 41' + [char]11 + '  4    1' + [char]11 + '  ----' + [char]11 + '5|    |' + [char]11 + '5|    |'
$cell.Range.Text = $new

$cell = $t.Cell(4, 2)
$new = '81 x 73' + [char]11 + '  7    3' + [char]11 + '  ----' + [char]11 + '8|    |' + [char]11 + '1|    |'
$cell.Range.Text = $new

$cell = $t.Cell(4, 3)
$new = '43 x 27' + [char]11 + '  2    7' + [char]11 + '  ----' + [char]11 + '4|    |' + [char]11 + '3|    |'
$cell.Range.Text = $new

$cell = $t.Cell(5, 1)
$new = '48 x 36' + [char]11 + '  3    6' + [char]11 + '  ----' + [char]11 + '4|    |' + [char]11 + '8|    |'
$cell.Range.Text = $new

$cell = $t.Cell(5, 2)
$new = '53 x 96' + [char]11 + '  9    6' + [char]11 + '  ----' + [char]11 + '5|    |' + [char]11 + '3|    |'
$cell.Range.Text = $new

$cell = $t.Cell(5, 3)
$new = '74 x 78' + [char]11 + '  7    8' + [char]11 + '  ----' + [char]11 + '7|    |'
$cell.Range.Text = $new
